$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for the "netdecked" card list (A: Carta, B: Utilidad, C: Cantidad)
$data = @(
    @("Bonfire", "Starter", 3),
    @("Ash Blossom & Joyous Spring", "Defensive", 3),
    @("Fire King High Avatar Kirin", "Extender", 3),
    @("Snake-Eyes Poplar", "Starter", 2),
    @("Legendary Fire King Ponix", "Starter", 1),
    @("Sacred Fire King Garunix", "Extender", 1),
    @("Fire King Island", "Garnet", 1),
    @("Original Sinful Spoils - Snake-Eye", "Extender", 1),
    @("Fire King Sanctuary", "Combo piece", 1),
    @("Fire King Avatar Arvata", "Combo piece", 1),
    @("Maxx C", "Defensive", 3),
    @("Snake-Eye Ash", "Starter", 1),
    @("Called by the Grave", "Defensive", 2),
    @("Snake-Eyes Flamberge Dragon", "Garnet", 1),
    @("Diabellstar the Black Witch", "Starter", 3),
    @("WANTED: Seeket of Sinful Spoils", "Starter", 1),
    @("Fire King Sky Burn", "Defensive", 1),
    @("Crossout Designator", "Defensive", 1),
    @("Infinite Impermanence", "Defensive", 3),
    @("Effect Veiler", "Defensive", 3),
    @("Nibiru the Primal Being", "Defensive", 1),
    @("Snake-Eye Oak", "Extender", 1),
    @("Fire King Courtier Ulcanix", "Starter", 2),
    @("Triple Tactic Talent", "Non engine", 1),
    @("One for one", "Starter", 1),
    @("Sinful Spoils of Subversion - Snake-Eye", "Defensive", 1),
    @("Divine Temple of the Snale-Eye", "Garnet", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Columns.Item(1).ColumnWidth = 35.7109375

$ws.Range("G16").Select()
